$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 71
$ws.Range("C2").Value = 71

$ws.Range("B3").Value = 92
$ws.Range("C3").Value = 83

$ws.Range("B4").Value = 39
$ws.Range("C4").Value = 23

$ws.Range("B5").Value = 94
$ws.Range("C5").Value = 8

$ws.Range("B6").Value = 87
$ws.Range("C6").Value = 94

$ws.Range("B7").Value = 93
$ws.Range("C7").Value = 14

$ws.Range("B8").Value = 12
$ws.Range("C8").Value = 86

$ws.Range("B9").Value = 68
$ws.Range("C9").Value = 48

$ws.Range("B10").Value = 28
$ws.Range("C10").Value = 32

$ws.Range("B11").Value = 78
$ws.Range("C11").Value = 94
